$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Copy formatting (style) of A2 down to the new rows A10:A17
# so the new index cells inherit the bold/centered/bordered style (s="1").
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A10:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Step 2: Write cell values for rows 2-17, columns A-M (matching target data).

# Row 2 (index 0)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1.058298367644753
$ws.Range("C2").Value = 0.5
$ws.Range("D2").Value = 46317314548.46852
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = "linearization_heuristic"
$ws.Range("G2").Value = 0.3
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 42825582311.74051
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = "linearization_heuristic"
$ws.Range("L2").Value = 3
$ws.Range("M2").Value = 0

# Row 3 (index 1)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1.361992657362816
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = 52792010116.85809
$ws.Range("E3").Value = 0.1
$ws.Range("F3").Value = "linearization_heuristic"
$ws.Range("G3").Value = 0.3
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 48298273933.09697
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = "linearization_heuristic"
$ws.Range("L3").Value = 3
$ws.Range("M3").Value = 0

# Row 4 (index 2)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1.108517102734018
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 46441108725.47605
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = "linearization_heuristic"
$ws.Range("G4").Value = 0.3
$ws.Range("H4").Value = 14
$ws.Range("I4").Value = 42783685619.74567
$ws.Range("J4").Value = 7
$ws.Range("K4").Value = "linearization_heuristic"
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 0

# Row 5 (index 3)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 1.424606802960026
$ws.Range("C5").Value = 0.5
$ws.Range("D5").Value = 52917919417.56814
$ws.Range("E5").Value = 0.1
$ws.Range("F5").Value = "linearization_heuristic"
$ws.Range("G5").Value = 0.3
$ws.Range("H5").Value = 14
$ws.Range("I5").Value = 48217595151.14416
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = "linearization_heuristic"
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 0

# Row 6 (index 4)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1.058298367644753
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 46317314548.46852
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("G6").Value = 0.3
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 42825582311.74051
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = "linearization_heuristic"
$ws.Range("L6").Value = 3
$ws.Range("M6").Value = 0

# Row 7 (index 5)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 1.361992657362816
$ws.Range("C7").Value = 0.5
$ws.Range("D7").Value = 52792010116.85809
$ws.Range("E7").Value = 0.1
$ws.Range("F7").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("G7").Value = 0.3
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 48298273933.09697
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = "linearization_heuristic"
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 0

# Row 8 (index 6)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 1.108517102734018
$ws.Range("C8").Value = 0.5
$ws.Range("D8").Value = 46441108725.47605
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("G8").Value = 0.3
$ws.Range("H8").Value = 14
$ws.Range("I8").Value = 42783685619.74567
$ws.Range("J8").Value = 7
$ws.Range("K8").Value = "linearization_heuristic"
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 0

# Row 9 (index 7)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 1.424606802960026
$ws.Range("C9").Value = 0.5
$ws.Range("D9").Value = 52917919417.56814
$ws.Range("E9").Value = 0.1
$ws.Range("F9").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("G9").Value = 0.3
$ws.Range("H9").Value = 14
$ws.Range("I9").Value = 48217595151.14416
$ws.Range("J9").Value = 7
$ws.Range("K9").Value = "linearization_heuristic"
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = 0

# Row 10 (index 8)
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 0.6034344150864958
$ws.Range("C10").Value = 0.5
$ws.Range("D10").Value = 43147341045.45089
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = "linearization_heuristic"
$ws.Range("G10").Value = 0.3
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 34422226983.28634
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = "linearization_heuristic"
$ws.Range("L10").Value = 3
$ws.Range("M10").Value = 11159709000

# Row 11 (index 9)
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 0.6271818492969048
$ws.Range("C11").Value = 0.5
$ws.Range("D11").Value = 48569441960.47471
$ws.Range("E11").Value = 0.1
$ws.Range("F11").Value = "linearization_heuristic"
$ws.Range("G11").Value = 0.3
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 39500961548.9781
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = "linearization_heuristic"
$ws.Range("L11").Value = 3
$ws.Range("M11").Value = 11159709000

# Row 12 (index 10)
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 0.6021885435055344
$ws.Range("C12").Value = 0.5
$ws.Range("D12").Value = 42846866826.35518
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = "linearization_heuristic"
$ws.Range("G12").Value = 0.3
$ws.Range("H12").Value = 14
$ws.Range("I12").Value = 34139766936.69876
$ws.Range("J12").Value = 7
$ws.Range("K12").Value = "linearization_heuristic"
$ws.Range("L12").Value = 3
$ws.Range("M12").Value = 11159709000

# Row 13 (index 11)
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 0.6448691365100838
$ws.Range("C13").Value = 0.5
$ws.Range("D13").Value = 48547780458.98592
$ws.Range("E13").Value = 0.1
$ws.Range("F13").Value = "linearization_heuristic"
$ws.Range("G13").Value = 0.3
$ws.Range("H13").Value = 14
$ws.Range("I13").Value = 39223557924.52286
$ws.Range("J13").Value = 7
$ws.Range("K13").Value = "linearization_heuristic"
$ws.Range("L13").Value = 3
$ws.Range("M13").Value = 11159709000

# Row 14 (index 12)
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 0.6034344150864958
$ws.Range("C14").Value = 0.5
$ws.Range("D14").Value = 43147341045.45089
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("G14").Value = 0.3
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 34422226983.28634
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = "linearization_heuristic"
$ws.Range("L14").Value = 3
$ws.Range("M14").Value = 11159709000

# Row 15 (index 13)
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 0.6271818492969047
$ws.Range("C15").Value = 0.5
$ws.Range("D15").Value = 48569441960.47471
$ws.Range("E15").Value = 0.1
$ws.Range("F15").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("G15").Value = 0.3
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 39500961548.9781
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = "linearization_heuristic"
$ws.Range("L15").Value = 3
$ws.Range("M15").Value = 11159709000

# Row 16 (index 14)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 0.6021885435055344
$ws.Range("C16").Value = 0.5
$ws.Range("D16").Value = 42846866826.35518
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("G16").Value = 0.3
$ws.Range("H16").Value = 14
$ws.Range("I16").Value = 34139766936.69876
$ws.Range("J16").Value = 7
$ws.Range("K16").Value = "linearization_heuristic"
$ws.Range("L16").Value = 3
$ws.Range("M16").Value = 11159709000

# Row 17 (index 15)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 0.6448691365100838
$ws.Range("C17").Value = 0.5
$ws.Range("D17").Value = 48547780458.98592
$ws.Range("E17").Value = 0.1
$ws.Range("F17").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("G17").Value = 0.3
$ws.Range("H17").Value = 14
$ws.Range("I17").Value = 39223557924.52286
$ws.Range("J17").Value = 7
$ws.Range("K17").Value = "linearization_heuristic"
$ws.Range("L17").Value = 3
$ws.Range("M17").Value = 11159709000
